# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" sheet with the latest case counts.
# Because the source rows are kept sorted by total cases (column B,
# descending), several countries swap row positions as their counts
# change relative to their neighbours. Each affected row below is
# rewritten in full (country name + all 7 numeric columns) so both
# the reordering and the updated figures land correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 31 de Marzo de 2020 a las 12:50
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 31 de Marzo de 2020 a las 12:50"

# Row 23: Noruega
$ws.Cells.Item(23, 1).Value = "Noruega"
$ws.Cells.Item(23, 2).Value = 4495
$ws.Cells.Item(23, 3).Value = 50
$ws.Cells.Item(23, 4).Value = 13
$ws.Cells.Item(23, 5).Value = 4448
$ws.Cells.Item(23, 6).Value = 97
$ws.Cells.Item(23, 7).Value = 2
$ws.Cells.Item(23, 8).Value = 34

# Row 25: Chequia
$ws.Cells.Item(25, 1).Value = "Chequia"
$ws.Cells.Item(25, 2).Value = 3002
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 25
$ws.Cells.Item(25, 5).Value = 2952
$ws.Cells.Item(25, 6).Value = 64
$ws.Cells.Item(25, 7).Value = 2
$ws.Cells.Item(25, 8).Value = 25

# Row 32: Polonia
$ws.Cells.Item(32, 1).Value = "Polonia"
$ws.Cells.Item(32, 2).Value = 2132
$ws.Cells.Item(32, 3).Value = 77
$ws.Cells.Item(32, 4).Value = 7
$ws.Cells.Item(32, 5).Value = 2094
$ws.Cells.Item(32, 6).Value = 50
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 31

# Row 95: Senegal
$ws.Cells.Item(95, 1).Value = "Senegal"
$ws.Cells.Item(95, 2).Value = 175
$ws.Cells.Item(95, 3).Value = 13
$ws.Cells.Item(95, 4).Value = 40
$ws.Cells.Item(95, 5).Value = 135
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 0

# Row 96: Afganistan
$ws.Cells.Item(96, 1).Value = "Afganistan"
$ws.Cells.Item(96, 2).Value = 174
$ws.Cells.Item(96, 3).Value = 4
$ws.Cells.Item(96, 4).Value = 5
$ws.Cells.Item(96, 5).Value = 165
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 4

# Row 97: Cuba
$ws.Cells.Item(97, 1).Value = "Cuba"
$ws.Cells.Item(97, 2).Value = 170
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = 4
$ws.Cells.Item(97, 5).Value = 162
$ws.Cells.Item(97, 6).Value = 2
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 4

# Row 98: Malta
$ws.Cells.Item(98, 1).Value = "Malta"
$ws.Cells.Item(98, 2).Value = 169
$ws.Cells.Item(98, 3).Value = 13
$ws.Cells.Item(98, 4).Value = 2
$ws.Cells.Item(98, 5).Value = 167
$ws.Cells.Item(98, 6).Value = 4
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0

# Row 99: Islas Feroe
$ws.Cells.Item(99, 1).Value = "Islas Feroe"
$ws.Cells.Item(99, 2).Value = 169
$ws.Cells.Item(99, 3).Value = 1
$ws.Cells.Item(99, 4).Value = 74
$ws.Cells.Item(99, 5).Value = 95
$ws.Cells.Item(99, 6).Value = 3
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0

# Row 100: Costa de Marfil
$ws.Cells.Item(100, 1).Value = "Costa de Marfil"
$ws.Cells.Item(100, 2).Value = 168
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 4).Value = 6
$ws.Cells.Item(100, 5).Value = 161
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 1

# Row 101: Uzbekistan
$ws.Cells.Item(101, 1).Value = "Uzbekistan"
$ws.Cells.Item(101, 2).Value = 158
$ws.Cells.Item(101, 3).Value = 9
$ws.Cells.Item(101, 4).Value = 7
$ws.Cells.Item(101, 5).Value = 149
$ws.Cells.Item(101, 6).Value = 8
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 2

# Row 106: Nigeria
$ws.Cells.Item(106, 1).Value = "Nigeria"
$ws.Cells.Item(106, 2).Value = 135
$ws.Cells.Item(106, 3).Value = 4
$ws.Cells.Item(106, 4).Value = 8
$ws.Cells.Item(106, 5).Value = 125
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 2

# Row 107: Venezuela
$ws.Cells.Item(107, 1).Value = "Venezuela"
$ws.Cells.Item(107, 2).Value = 135
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(107, 4).Value = 39
$ws.Cells.Item(107, 5).Value = 93
$ws.Cells.Item(107, 6).Value = 6
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 3

# Row 143: Niger
$ws.Cells.Item(143, 1).Value = "Niger"
$ws.Cells.Item(143, 2).Value = 27
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 24
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 3

# Row 144: Bermudas
$ws.Cells.Item(144, 1).Value = "Bermudas"
$ws.Cells.Item(144, 2).Value = 27
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 10
$ws.Cells.Item(144, 5).Value = 17
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 0

# Row 145: Etiopia
$ws.Cells.Item(145, 1).Value = "Etiopia"
$ws.Cells.Item(145, 2).Value = 25
$ws.Cells.Item(145, 3).Value = 2
$ws.Cells.Item(145, 4).Value = 2
$ws.Cells.Item(145, 5).Value = 23
$ws.Cells.Item(145, 6).Value = 2
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 0

# Row 146: Mali
$ws.Cells.Item(146, 1).Value = "Mali"
$ws.Cells.Item(146, 2).Value = 25
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 0
$ws.Cells.Item(146, 5).Value = 23
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 2

# Row 167: Groenlandia
$ws.Cells.Item(167, 1).Value = "Groenlandia"
$ws.Cells.Item(167, 2).Value = 10
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 2
$ws.Cells.Item(167, 5).Value = 8
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 0

# Row 168: Siria
$ws.Cells.Item(168, 1).Value = "Siria"
$ws.Cells.Item(168, 2).Value = 10
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 0
$ws.Cells.Item(168, 5).Value = 8
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 2

# Row 173: Libia
$ws.Cells.Item(173, 1).Value = "Libia"
$ws.Cells.Item(173, 2).Value = 8
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 8
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

# Row 174: Mozambique
$ws.Cells.Item(174, 1).Value = "Mozambique"
$ws.Cells.Item(174, 2).Value = 8
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 8
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

# Row 175: Guinea-Bisau
$ws.Cells.Item(175, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(175, 2).Value = 8
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = 8
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

# Row 176: Surinam
$ws.Cells.Item(176, 1).Value = "Surinam"
$ws.Cells.Item(176, 2).Value = 8
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 5).Value = 8
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

# Row 177: Guyana
$ws.Cells.Item(177, 1).Value = "Guyana"
$ws.Cells.Item(177, 2).Value = 8
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 7
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 1

# Row 178: Zimbabue
$ws.Cells.Item(178, 1).Value = "Zimbabue"
$ws.Cells.Item(178, 2).Value = 8
$ws.Cells.Item(178, 3).Value = 1
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 5).Value = 7
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 1

# Row 184: Santa Sede
$ws.Cells.Item(184, 1).Value = "Santa Sede"
$ws.Cells.Item(184, 2).Value = 6
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 6
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

# Row 185: San Martin (Parte Holandesa)
$ws.Cells.Item(185, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(185, 2).Value = 6
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 6
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

# Row 186: Cabo Verde
$ws.Cells.Item(186, 1).Value = "Cabo Verde"
$ws.Cells.Item(186, 2).Value = 6
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 5
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 1

# Row 187: San Bartolome
$ws.Cells.Item(187, 1).Value = "San Bartolome"
$ws.Cells.Item(187, 2).Value = 6
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 1
$ws.Cells.Item(187, 5).Value = 5
$ws.Cells.Item(187, 6).Value = 0
$ws.Cells.Item(187, 7).Value = 0
$ws.Cells.Item(187, 8).Value = 0

# Row 188: Benin
$ws.Cells.Item(188, 1).Value = "Benin"
$ws.Cells.Item(188, 2).Value = 6
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 1
$ws.Cells.Item(188, 5).Value = 5
$ws.Cells.Item(188, 6).Value = 0
$ws.Cells.Item(188, 7).Value = 0
$ws.Cells.Item(188, 8).Value = 0

# Row 190: Islas Turcas y Caicos
$ws.Cells.Item(190, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(190, 2).Value = 5
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 5).Value = 5
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

# Row 192: Montserrat
$ws.Cells.Item(192, 1).Value = "Montserrat"
$ws.Cells.Item(192, 2).Value = 5
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

# Row 195: Gambia
$ws.Cells.Item(195, 1).Value = "Gambia"
$ws.Cells.Item(195, 2).Value = 4
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 3
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 1

# Row 196: Nicaragua
$ws.Cells.Item(196, 1).Value = "Nicaragua"
$ws.Cells.Item(196, 2).Value = 4
$ws.Cells.Item(196, 3).Value = 0
$ws.Cells.Item(196, 4).Value = 0
$ws.Cells.Item(196, 5).Value = 3
$ws.Cells.Item(196, 6).Value = 0
$ws.Cells.Item(196, 7).Value = 0
$ws.Cells.Item(196, 8).Value = 1

# Row 197: Belice
$ws.Cells.Item(197, 1).Value = "Belice"
$ws.Cells.Item(197, 2).Value = 3
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 5).Value = 3
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 0

# Row 198: Republica de Africa Central
$ws.Cells.Item(198, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(198, 2).Value = 3
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 0
$ws.Cells.Item(198, 5).Value = 3
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 0

# Row 204: Papua Nueva Guinea
$ws.Cells.Item(204, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(204, 2).Value = 1
$ws.Cells.Item(204, 3).Value = 0
$ws.Cells.Item(204, 4).Value = 0
$ws.Cells.Item(204, 5).Value = 1
$ws.Cells.Item(204, 6).Value = 0
$ws.Cells.Item(204, 7).Value = 0
$ws.Cells.Item(204, 8).Value = 0

# Row 205: Timor Oriental
$ws.Cells.Item(205, 1).Value = "Timor Oriental"
$ws.Cells.Item(205, 2).Value = 1
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 0
$ws.Cells.Item(205, 5).Value = 1
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0
